$wb = $excel.ActiveWorkbook

# --- Personnel sheet edits ---
$personnel = $wb.Worksheets.Item("Personnel")

# Correct Kevin Cahill's role from "technician" to "creator"
$personnel.Range("G6").Value = "creator"

# Insert a new row for Zoe Sandwith (AR28) as a creator, pushing Kate Morkeski down to row 8
$personnel.Rows.Item(7).Insert()
$personnel.Range("A7").Value = "Zoe"
$personnel.Range("C7").Value = "Sandwith"
$personnel.Range("D7").Value = "Northeast U.S. Shelf LTER"
$personnel.Range("G7").Value = "creator"

# --- Sheet view / selection updates ---
$categorical = $wb.Worksheets.Item("CategoricalVariables")
$categorical.Range("A4").Select()

$personnel.Range("E11").Select()
$personnel.Activate()
